# PLD.xlsx - MarketBeatRank weekly refresh ("10Th - MB for single stock and
# added new group"):
#   1. Two new date snapshots (Jun_26, Jun_27) are pulled in as three new
#      leading data columns (B, C, D) - the existing date columns (old
#      B:E = Jun_17/Jun_15/Jun_13/Jun_10) shift right to become E:H.
#   2. Two brokerages that just initiated coverage ("Benchmark" and
#      "Evercore ISI") are appended as new rows 28/29, populated only for
#      the three newest snapshot columns (A:D) since they have no history
#      for the older dates.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Make room for the 3 newest snapshots -------------------------------
# Inserting a 3-column-wide selection at B pushes old B:E -> E:H in one shot.
$ws.Range("B1:D1").EntireColumn.Insert()
$ws.Range("B:D").ColumnWidth = 41.67

# --- 2. Fill the two "Jun_26" columns (D then C) ----------------------------
$ws.Range("D1").Value = "Jun_26"
$ws.Range("D2:D27").Value = "UN"

$ws.Range("C1").Value = "Jun_26"
$ws.Range("C2:C27").Value = "UN"

# --- 3. Append the two newly-covering brokerages ----------------------------
$ws.Range("A28").Value = "Benchmark"
$ws.Range("B28:D28").Value = "UN"

$ws.Range("A29").Value = "Evercore ISI"
$ws.Range("B29:D29").Value = "UN"

# --- 4. Fill the newest "Jun_27" column (B) ---------------------------------
$ws.Range("B1").Value = "Jun_27"
$ws.Range("B2:B27").Value = "UN"
